$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the old "Voice ID" placeholder letters (A / B) with real voice IDs,
# and add a new "Voice Name" column next to it.
$ws.Range("D1").Value = "Voice Name"

$ws.Range("C2").Value = "t0jbNlBVZ17f02VDIeMI"
$ws.Range("D2").Value = "Jessie"

$ws.Range("C3").Value = "zcAOhNBS3c14rBihAFp1"
$ws.Range("D3").Value = "Giovanni"

# Auto-fit the Voice ID column now that it holds longer values, matching the
# author's "best fit" width the way column B already got.
$ws.Columns.Item(3).EntireColumn.AutoFit() | Out-Null
$ws.Columns.Item(3).ColumnWidth = 18.25

# Leave the selection where the author ended up after entering the data.
$ws.Range("F9").Select()
